# Scheduled-runner update: refresh market-price-derived profit figures
# across the per-job Leve sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Only numeric value cells change; no structural edits.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 51013.4
$ws.Range("I28").Value = 67183.2
$ws.Range("K28").Value = 67183.2
$ws.Range("M28").Value = -66698.2
$ws.Range("H64").Value = 6074.077
$ws.Range("I64").Value = 4999.8
$ws.Range("K64").Value = 4999.8
$ws.Range("M64").Value = -4751.8
$ws.Range("H67").Value = 6074.077
$ws.Range("I67").Value = 4999.8
$ws.Range("K67").Value = 4999.8
$ws.Range("M67").Value = -4141.8
$ws.Range("H113").Value = 6904.25
$ws.Range("I113").Value = 3500
$ws.Range("J113").Value = 7585.1
$ws.Range("K113").Value = 3500
$ws.Range("L113").Value = 7585.1
$ws.Range("M113").Value = -246
$ws.Range("N113").Value = -14093.1
$ws.Range("H135").Value = 1002.1667
$ws.Range("I135").Value = 852.5333000000001
$ws.Range("J135").Value = 1750.3334
$ws.Range("K135").Value = 7672.7997
$ws.Range("L135").Value = 15753.0006
$ws.Range("M135").Value = -5137.7997
$ws.Range("N135").Value = -20823.0006
$ws.Range("H137").Value = 3152.8262
$ws.Range("I137").Value = 2111.75
$ws.Range("J137").Value = 3708.0667
$ws.Range("K137").Value = 6335.25
$ws.Range("L137").Value = 11124.2001
$ws.Range("M137").Value = -3785.25
$ws.Range("N137").Value = -16224.2001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 1668.3334
$ws.Range("J3").Value = 4006
$ws.Range("L3").Value = 4006
$ws.Range("N3").Value = -4236
$ws.Range("H8").Value = 7000
$ws.Range("I8").Value = 7000
$ws.Range("K8").Value = 7000
$ws.Range("M8").Value = -6856
$ws.Range("H11").Value = 1400
$ws.Range("I11").Value = 1300
$ws.Range("K11").Value = 1300
$ws.Range("M11").Value = -1156
$ws.Range("H13").Value = 6274.5
$ws.Range("I13").Value = 1500
$ws.Range("J13").Value = 7866
$ws.Range("K13").Value = 1500
$ws.Range("L13").Value = 7866
$ws.Range("M13").Value = -1356
$ws.Range("N13").Value = -8154
$ws.Range("H74").Value = 25644276
$ws.Range("J74").Value = 2400
$ws.Range("L74").Value = 2400
$ws.Range("N74").Value = -4148
$ws.Range("H77").Value = 25644276
$ws.Range("J77").Value = 2400
$ws.Range("L77").Value = 12000
$ws.Range("N77").Value = -20736
$ws.Range("H132").Value = 3065.9285
$ws.Range("I132").Value = 2479.6
$ws.Range("J132").Value = 4531.75
$ws.Range("K132").Value = 7438.799999999999
$ws.Range("L132").Value = 13595.25
$ws.Range("M132").Value = -4908.799999999999
$ws.Range("N132").Value = -18655.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1821.8
$ws.Range("I5").Value = 118.25
$ws.Range("K5").Value = 118.25
$ws.Range("M5").Value = -5.25
$ws.Range("H134").Value = 1818.091
$ws.Range("I134").Value = 1175.75
$ws.Range("K134").Value = 3527.25
$ws.Range("M134").Value = -992.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 44073.777
$ws.Range("J31").Value = 86464.16
$ws.Range("L31").Value = 86464.16
$ws.Range("N31").Value = -87054.16
$ws.Range("H32").Value = 950.5
$ws.Range("I32").Value = 950.5
$ws.Range("K32").Value = 950.5
$ws.Range("M32").Value = -634.5
$ws.Range("H34").Value = 44073.777
$ws.Range("J34").Value = 86464.16
$ws.Range("L34").Value = 86464.16
$ws.Range("N34").Value = -86868.16
$ws.Range("H41").Value = 20000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 20000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 20000
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -20856
$ws.Range("H107").Value = 916.2692
$ws.Range("J107").Value = 1175.125
$ws.Range("L107").Value = 1175.125
$ws.Range("N107").Value = -5015.125
$ws.Range("H132").Value = 11445.777
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 53.625
$ws.Range("I12").Value = 67.333336
$ws.Range("J12").Value = 45.4
$ws.Range("K12").Value = 202.000008
$ws.Range("L12").Value = 136.2
$ws.Range("M12").Value = -29.00000800000001
$ws.Range("N12").Value = -482.2
$ws.Range("H86").Value = 2129.8
$ws.Range("J86").Value = 3232.6667
$ws.Range("L86").Value = 9698.000100000001
$ws.Range("N86").Value = -12070.0001
$ws.Range("H89").Value = 2129.8
$ws.Range("J89").Value = 3232.6667
$ws.Range("L89").Value = 29094.0003
$ws.Range("N89").Value = -40950.0003
$ws.Range("H102").Value = 10000
$ws.Range("J102").Value = 10000
$ws.Range("L102").Value = 30000
$ws.Range("N102").Value = -34868
$ws.Range("H140").Value = 3149.75
$ws.Range("I140").Value = 3456.8572
$ws.Range("K140").Value = 10370.5716
$ws.Range("M140").Value = -5190.571599999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16730.422
$ws.Range("I70").Value = 6277.357
$ws.Range("K70").Value = 6277.357
$ws.Range("M70").Value = -6007.357
$ws.Range("H73").Value = 16730.422
$ws.Range("I73").Value = 6277.357
$ws.Range("K73").Value = 6277.357
$ws.Range("M73").Value = -5341.357
$ws.Range("H80").Value = 6301.857
$ws.Range("I80").Value = 4765.8335
$ws.Range("J80").Value = 6916.2666
$ws.Range("K80").Value = 4765.8335
$ws.Range("L80").Value = 6916.2666
$ws.Range("M80").Value = -3767.8335
$ws.Range("N80").Value = -8912.266599999999
$ws.Range("H83").Value = 6301.857
$ws.Range("I83").Value = 4765.8335
$ws.Range("J83").Value = 6916.2666
$ws.Range("K83").Value = 23829.1675
$ws.Range("L83").Value = 34581.333
$ws.Range("M83").Value = -18837.1675
$ws.Range("N83").Value = -44565.333

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 4868
$ws.Range("I12").Value = 10000
$ws.Range("J12").Value = 2302
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 2302
$ws.Range("M12").Value = -9830
$ws.Range("N12").Value = -2642
$ws.Range("H22").Value = 8300.5
$ws.Range("I22").Value = 2100
$ws.Range("K22").Value = 2100
$ws.Range("M22").Value = -1805
$ws.Range("H27").Value = 8300.5
$ws.Range("I27").Value = 2100
$ws.Range("K27").Value = 2100
$ws.Range("M27").Value = -1993
$ws.Range("H68").Value = 3694.348
$ws.Range("I68").Value = 2998.5
$ws.Range("K68").Value = 2998.5
$ws.Range("M68").Value = -2249.5
$ws.Range("H71").Value = 3694.348
$ws.Range("I71").Value = 2998.5
$ws.Range("K71").Value = 14992.5
$ws.Range("M71").Value = -11248.5
$ws.Range("H82").Value = 6344.923
$ws.Range("I82").Value = 6373.625
$ws.Range("J82").Value = 6299
$ws.Range("K82").Value = 6373.625
$ws.Range("L82").Value = 6299
$ws.Range("M82").Value = -6012.625
$ws.Range("N82").Value = -7021
$ws.Range("H85").Value = 6344.923
$ws.Range("I85").Value = 6373.625
$ws.Range("J85").Value = 6299
$ws.Range("K85").Value = 6373.625
$ws.Range("L85").Value = 6299
$ws.Range("M85").Value = -5125.625
$ws.Range("N85").Value = -8795
$ws.Range("H132").Value = 4192.6206
$ws.Range("I132").Value = 2584.85
$ws.Range("K132").Value = 7754.549999999999
$ws.Range("M132").Value = -5224.549999999999
$ws.Range("H136").Value = 13651.429
$ws.Range("I136").Value = 3750
$ws.Range("J136").Value = 15301.667
$ws.Range("K136").Value = 11250
$ws.Range("L136").Value = 45905.001
$ws.Range("M136").Value = -8700
$ws.Range("N136").Value = -51005.001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H62").Value = 8487.875
$ws.Range("I62").Value = 8301.333000000001
$ws.Range("K62").Value = 8301.333000000001
$ws.Range("M62").Value = -7677.333000000001
$ws.Range("H65").Value = 8487.875
$ws.Range("I65").Value = 8301.333000000001
$ws.Range("K65").Value = 41506.665
$ws.Range("M65").Value = -38386.665
$ws.Range("H132").Value = 5786.087
$ws.Range("I132").Value = 5056.8945
$ws.Range("K132").Value = 15170.6835
$ws.Range("M132").Value = -12640.6835

